$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.623.84'
$ws.Range('E2').Value = '  +1.22%  '

$ws.Range('D3').Value = '2.986.79'
$ws.Range('E3').Value = '  +2.73%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = "'384.62"
$ws.Range('E5').Value = '  +2.83%  '

$ws.Range('D6').Value = "'104.64"
$ws.Range('E6').Value = '  +3.14%  '

$ws.Range('E7').Value = '  +0.96%  '

$ws.Range('E8').Value = '  +0.15%  '

$ws.Range('D9').Value = "'0.596"
$ws.Range('E9').Value = '  +1.98%  '

$ws.Range('D10').Value = "'37.24"
$ws.Range('E10').Value = '  +1.38%  '

$ws.Range('D12').Value = "'0.0850"
$ws.Range('E12').Value = '  +2.04%  '

$ws.Range('D13').Value = '3.461.26'
$ws.Range('E13').Value = '  +2.73%  '

$ws.Range('D14').Value = "'18.40"
$ws.Range('E14').Value = '  +0.98%  '

$ws.Range('D15').Value = "'7.61"
$ws.Range('E15').Value = '  +3.63%  '

$ws.Range('D16').Value = '2.988.03'
$ws.Range('E16').Value = '  +2.82%  '

$ws.Range('E17').Value = '  +9.37%  '

$ws.Range('D18').Value = '51.617.30'
$ws.Range('E18').Value = '  +1.26%  '

$ws.Range('D19').Value = "'3.28"
$ws.Range('E19').Value = '  +1.40%  '

$ws.Range('D20').Value = "'7.46"
$ws.Range('E20').Value = '  +3.70%  '

$ws.Range('D21').Value = "'12.92"
$ws.Range('E21').Value = '  +0.97%  '

$ws.Range('D22').Value = '0.0₃0965'
$ws.Range('E22').Value = '  +2.55%  '

$ws.Range('D23').Value = "'69.17"
$ws.Range('E23').Value = '  +1.51%  '

$ws.Range('D24').Value = "'263.61"
$ws.Range('E24').Value = '  +1.77%  '

$ws.Range('E25').Value = '  +8.81%  '

$ws.Range('D26').Value = "'8.36"
$ws.Range('E26').Value = '  +18.91%  '

$ws.Range('D27').Value = "'7.72"
$ws.Range('E27').Value = '  +21.80%  '

$ws.Range('B28').Value = 'Hedera'
$ws.Range('C28').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D28').Value = "'0.115"
$ws.Range('E28').Value = '  +14.76%  '

$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').Value = "'0.170"
$ws.Range('E29').Value = '  +1.31%  '

$ws.Range('D30').Value = "'26.05"

$ws.Range('E31').Value = '  -0.11%  '

$ws.Range('D32').Value = "'9.91"
$ws.Range('E32').Value = '  +0.96%  '

$ws.Range('D33').Value = "'34.89"
$ws.Range('E33').Value = '  +2.51%  '

$ws.Range('D34').Value = "'51.06"
$ws.Range('E34').Value = '  -0.32%  '

$ws.Range('E35').Value = '  -1.70%  '

$ws.Range('D36').Value = "'0.0453"
$ws.Range('E36').Value = '  +7.66%  '

$ws.Range('E37').Value = '  -0.03%  '

$ws.Range('E38').Value = '  +2.17%  '

$ws.Range('D39').Value = "'17.06"
$ws.Range('E39').Value = '  +0.87%  '

$ws.Range('E40').Value = '  +1.54%  '

$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').Value = "'0.116"
$ws.Range('E41').Value = '  +3.72%  '

$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = "'1.84"
$ws.Range('E42').Value = '  +0.58%  '

$ws.Range('D43').Value = "'122.40"
$ws.Range('E43').Value = '  +2.81%  '

$ws.Range('D44').Value = "'21.78"
$ws.Range('E44').Value = '  -0.17%  '

$ws.Range('E45').Value = '  +18.36%  '

$ws.Range('E46').Value = '  -1.81%  '

$ws.Range('E47').Value = '  +2.74%  '

$ws.Range('E48').Value = '  +5.44%  '

$ws.Range('D49').Value = '2.038.56'
$ws.Range('E49').Value = '  +1.36%  '

$ws.Range('D50').Value = "'0.0333"
$ws.Range('E50').Value = '  +8.42%  '

$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = "'5.16"
$ws.Range('E51').Value = '  +2.74%  '

